# "product page not finish"
# - the green "done" highlight is removed from the whole progress table,
#   leaving every row highlighted yellow (still in progress)
# - row 10 (#5, previously just a bare index with no task filled in yet)
#   gets its task name / sub-task filled in
# - a new row 11 (#6) is appended for the next not-yet-finished task

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535   # OLE BGR encoding of RGB(255,255,0) - matches the sheet's existing yellow fill

# 1) Strip the green fill from every previously-coloured cell (rows 3-9),
#    leaving the yellow fill that was already used elsewhere - alignment
#    (center / center+middle) for each cell is left untouched.
$ws.Range("A3:C9").Interior.Color = $yellow

# 2) Row 10 ("5") only had its index filled in - give it the same yellow,
#    centered look as the rest of the table and fill in the task text.
$ws.Range("A10:C10").Interior.Color = $yellow
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("B10").Value = "設計產品畫面"
$ws.Range("C10").Value = "添加cart button"
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("C10").HorizontalAlignment = -4108

# 3) New row 11 ("6") for the next outstanding task - plain (not yet
#    highlighted), matching the sheet's default column style.
$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "設計個別產品畫面"

# 4) Leave the cursor where the author left it when they saved.
$ws.Range("C11").Select()
